$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "sodium"
$ws.Range("B2").Value = "ORG"
$ws.Range("C2").Value = 658
$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.wikidata.org/wiki/Q658")
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = 3
